$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column A width: 13.4 -> 23.48 chars (closest achievable via COM rounding)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.6

# ---------------------------------------------------------------------------
# 2. Row 4 new cell H4: shipstats value for the colonists-overheat event
# ---------------------------------------------------------------------------
$ws.Range("H4").Value = 'max_systems=1, max_damage=75,damageable="colonists"'

# ---------------------------------------------------------------------------
# 3. New events rows 5-7 (scanner / gearlanding / dbase overheat variants)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "events/overheat/scanner"
$ws.Range("B5").Value = "Allow the scanner module to overheat"
$ws.Range("E5").Value = "The AI channels excess heat into the scanner array. The heat regulation system recovers as the ship moves away from the super-hot star, but not before the heat has damaged a scanner."
$ws.Range("H5").Value = 'max_systems=1, max_damage=15,damageable=c("resources_sensor","temperature_sensor", "gravity_sensor", "atmosphere_sensor","water_sensor")'
$ws.Range("L5").Value = 0

$ws.Range("A6").Value = "events/overheat/gearlanding"
$ws.Range("B6").Value = "Allow the landing/construction module to overheat"
$ws.Range("E6").Value = "The AI channels excess heat into the landing gear and construction equipment. The heat regulation system recovers as the ship moves away from the super-hot star, but not before some heat damage is done."
$ws.Range("H6").Value = 'max_systems=1, max_damage=2,damageable=c("landing_gear","equipment")'
$ws.Range("L6").Value = 0

$ws.Range("A7").Value = "events/overheat/dbase"
$ws.Range("B7").Value = "Allow the data storage module to overheat"
$ws.Range("E7").Value = "The AI channels excess heat into the data storage module. The heat regulation system recovers as the ship moves away from the super-hot star, but not before some data is lost."
$ws.Range("H7").Value = 'max_systems=1, max_damage=7,damageable=c("planetLocalDB","dbase")'
$ws.Range("L7").Value = 0

# ---------------------------------------------------------------------------
# 4. Formatting: wrap text on the long description column for the new rows
#    (matches the wrapped "description" style used on rows 3-4, column E)
# ---------------------------------------------------------------------------
$ws.Range("E5:E7").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Row heights for the new rows
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 82.9
$ws.Rows.Item(6).RowHeight = 94.5
$ws.Rows.Item(7).RowHeight = 82.9

# ---------------------------------------------------------------------------
# 6. Touch two previously-empty cells so they materialize in the sheet (F3/F4)
# ---------------------------------------------------------------------------
$ws.Range("F3").NumberFormat = "General"
$ws.Range("F4").NumberFormat = "General"

# ---------------------------------------------------------------------------
# 7. View / selection state to match final cursor location after the edit
# ---------------------------------------------------------------------------
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("H7").Select()

Write-Host "done"
